$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A3").Value = "roll"
$ws.Range("B3").Value = "software"

$ws.Range("B3").Select()
